# Clean up vaccine/brand labels across every sheet:
#   - drop footnote markers like " [5]" (bracketed numbers)
#   - collapse embedded line breaks inside a cell into a single space
# This naturally de-duplicates the stray "Afluria\nQuadrivalent" shared
# string (Adult Influenza sheet, B9:B10) into the already-present
# "Afluria Quadrivalent" entry once the newline is gone.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ur = $ws.UsedRange
    $rowCount = $ur.Rows.Count
    $colCount = $ur.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ur.Cells.Item($r, $c)
            $txt = $cell.Text

            if ($txt -eq $null -or $txt -eq "") {
                continue
            }

            $new = $txt -replace '\[\d+\]', ''
            $new = $new -replace "`r`n", " "
            $new = $new -replace "`n", " "

            if ($new -ne $txt) {
                $cell.Value = $new
            }
        }
    }
}
